$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reel/coil rows appended below the existing table (rows 33-43).
# Force text formatting first so numeric-looking values (e.g. "04", "110")
# are stored as text like the rest of the sheet, then clear the temporary
# formatting so the cells keep the workbook's default (unstyled) look.
$newRange = $ws.Range("A33:K43")
$newRange.NumberFormat = "@"

$newData = @(
    ,@('62,5', '110', '347', '406', '14681', '1', '86989', '2025-08-22 13:16', 'B', '04', 'CART.GRIS')
    ,@('62,5', '110', '347', '406', '14681', '2', '86989', '2025-08-22 13:16', 'B', '04', 'CART.GRIS')
    ,@('62,5', '110', '347', '406', '14681', '3', '86989', '2025-08-22 13:16', 'B', '04', 'CART.GRIS')
    ,@('62,5', '110', '347', '406', '14680', '1', '86989', '2025-08-22 13:21', 'B', '04', 'CART.GRIS')
    ,@('62', '110', '347', '406', '14680', '1', '86989', '2025-08-22 13:23', 'B', '04', 'CART.GRIS')
    ,@('62,5', '110', '347', '406', '14680', '1', '86989', '2025-08-22 13:27', 'B', '04', 'CART.GRIS')
    ,@('62,5', '110', '347', '406', '14680', '1', '86989', '22/08/2025 13:32', 'B', '04', 'CART.GRIS')
    ,@('120', '120', '151', '120', '1544', '1', '45555', '25/08/2025 10:27', 'A', '03', 'L.BLANCO')
    ,@('120', '120', '151', '120', '1544', '1', '45555', '25/08/2025 10:27', 'A', '03', 'L.BLANCO')
    ,@('120', '120', '12', '1212', '1212', '1', '2221', '25/08/2025 10:31', 'A', '03', 'L.BLANCO')
    ,@('120', '120', '12', '1212', '1212', '2', '2221', '25/08/2025 10:31', 'A', '03', 'L.BLANCO')
)

$startRow = 33
for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowVals = $newData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

$newRange.ClearFormats()
